$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 (R1): uptime
$ws.Range("G3").Value = "0:07:00"

# Row 4 (R3): uptime
$ws.Range("G4").Value = "0:11:21"

# Row 5 (SW1): serial_number and uptime
$ws.Range("E5").Value = "9K3RW05NXRW"
$ws.Range("G5").Value = "0:11:00"
